$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card7")

# O1: remove trailing space -> "Serviced by"
$ws.Range("O1").Value = "Serviced by"

# L2: "nan" -> "1/8/2025" as literal text (not an Excel date).
# Force the cell to Text format before assigning so the date-like
# string isn't auto-parsed into a date serial number, then clear the
# formatting change back off so the cell style stays the default (0).
$L2 = $ws.Range("L2")
$L2.NumberFormat = "@"
$L2.Value = "1/8/2025"
$L2.ClearFormats()

# O2:O12 were empty inline-string cells; fill them with the text "nan"
# to match the rest of the row (mirrors columns D:N on each row).
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Value = "nan"
}
